# The "LoginData" sheet's A2 cell (hyperlinked e-mail used as a login
# test-fixture value) needs to be refreshed to a new address.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LoginData")
$ws.Range("A2").Value = "juan.perez99_46781@test.com"
